$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set A2 to the string value "s" (was numeric 0)
$ws.Range("A2").Value = "s"

# Remove row 3 entirely (was A3 = 456); this shrinks the used range to A1:A2
$ws.Rows.Item(3).Delete()
